$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '68.279.76'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.90%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.622.75'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.64%  '
$ws.Range('E4').Value = '  -0.36%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '583.15'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.52%  '
$ws.Range('E6').Value = '  +4.12%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.616.00'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.74%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.622'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.64%  '
$ws.Range('E9').Value = '  -0.03%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.681'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.26%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.155'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +4.67%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '55.80'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.56%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000294'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +12.98%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '10.10'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.68%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.184.43'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.43%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.615.53'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.25%  '
$ws.Range('E17').Value = '  -0.03%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.57'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.97%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '68.122.46'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.67%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '18.59'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.38%  '
$ws.Range('E21').Value = '  -0.50%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '405.19'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.24%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '13.06'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +22.98%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.25'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.12%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '86.21'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.13%  '
$ws.Range('B26').Value = 'ImmutableX'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.96'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.61%  '
$ws.Range('B27').Value = 'Toncoin'
$ws.Range('C27').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '4.00'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +10.58%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '12.64'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.98%  '
$ws.Range('E29').Value = '  +1.06%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.16'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +17.62%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '9.21'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.16%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '31.78'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.04%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '689.26'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +14.74%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '12.28'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.35%  '
$ws.Range('E35').Value = '  +4.48%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '64.88'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.34%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '42.85'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.04%  '
$ws.Range('E38').Value = '  +9.42%  '
$ws.Range('E39').Value = '  +0.11%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0₃0798'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +9.49%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.87'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +17.27%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.13'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +10.80%  '
$ws.Range('E43').Value = '  +0.39%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.145.01'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +15.35%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.998'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.41%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0425'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.67%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.133'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.29%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.91'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.40%  '
$ws.Range('B49').Value = 'Monero'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '144.02'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.31%  '
$ws.Range('B50').Value = 'dogwifhat'
$ws.Range('C50').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.64'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +11.49%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.63'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.39%  '
